$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 13:22"

# Rows 19-20: Malaga/Salamanca swap places (Malaga now ranks above Salamanca),
# each bringing its own refreshed totals.
$ws.Range("A19").Value = "Malaga"
$ws.Range("B19").Value = 2531
$ws.Range("C19").Value = 869
$ws.Range("D19").Value = 1439
$ws.Range("E19").Value = 223

$ws.Range("A20").Value = "Salamanca"
$ws.Range("B20").Value = 2514
$ws.Range("C20").Value = 761
$ws.Range("D20").Value = 1431
$ws.Range("E20").Value = 276

# Row 21: Sevilla keeps its place but gets updated figures.
$ws.Range("B21").Value = 2329
$ws.Range("C21").Value = 459
$ws.Range("D21").Value = 1658
$ws.Range("E21").Value = 212

# Rows 28-29: Granada/A Coruña swap places.
$ws.Range("A28").Value = "Granada"
$ws.Range("B28").Value = 2023
$ws.Range("C28").Value = 616
$ws.Range("D28").Value = 1211
$ws.Range("E28").Value = 196

$ws.Range("A29").Value = "A Coruña"
$ws.Range("B29").Value = 1969
$ws.Range("C29").Value = 333
$ws.Range("D29").Value = 1788
$ws.Range("E29").Value = 67

# Rows 35-36: Jaen/Cuenca swap places.
$ws.Range("A35").Value = "Jaen"
$ws.Range("B35").Value = 1297
$ws.Range("C35").Value = 253
$ws.Range("D35").Value = 907
$ws.Range("E35").Value = 137

$ws.Range("A36").Value = "Cuenca"
$ws.Range("B36").Value = 1285
$ws.Range("C36").Value = 3838
$ws.Range("D36").Value = 10545
$ws.Range("E36").Value = 154

# Row 37: Cordoba keeps its place but gets updated figures.
$ws.Range("B37").Value = 1271
$ws.Range("C37").Value = 371
$ws.Range("D37").Value = 826
$ws.Range("E37").Value = 74

# Row 40: Cadiz keeps its place but gets updated figures.
$ws.Range("B40").Value = 1139
$ws.Range("C40").Value = 283
$ws.Range("D40").Value = 781
$ws.Range("E40").Value = 75

# Row 51: Almeria keeps its place but gets updated figures.
$ws.Range("B51").Value = 458
$ws.Range("C51").Value = 147
$ws.Range("D51").Value = 269
$ws.Range("E51").Value = 42

# Row 52: Huelva keeps its place but gets updated figures.
$ws.Range("B52").Value = 377
$ws.Range("C52").Value = 118
$ws.Range("D52").Value = 225
$ws.Range("E52").Value = 34
